$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.036.02'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '2.307.71'
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '308.76'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').Value = '108.23'
$ws.Range('E6').Value = '  -6.00%  '
$ws.Range('D7').Value = '0.633'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').Value = '0.614'
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('D10').Value = '43.74'
$ws.Range('E10').Value = '  -5.34%  '
$ws.Range('D11').Value = '0.0927'
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('D12').Value = '8.77'
$ws.Range('E12').Value = '  -4.29%  '
$ws.Range('D13').Value = '1.04'
$ws.Range('E13').Value = '  +15.83%  '
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').Value = '15.73'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '2.649.22'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').Value = '2.344.29'
$ws.Range('E17').Value = '  +3.94%  '
$ws.Range('D18').Value = '43.083.12'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').Value = '7.22'
$ws.Range('E20').Value = '  -6.69%  '
$ws.Range('D21').Value = '75.80'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('D22').Value = '3.48'
$ws.Range('E22').Value = '  -6.66%  '
$ws.Range('E23').Value = '  +6.65%  '
$ws.Range('D24').Value = '255.34'
$ws.Range('E24').Value = '  +9.79%  '
$ws.Range('D25').Value = '8.95'
$ws.Range('E25').Value = '  -6.36%  '
$ws.Range('D26').Value = '11.83'
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('D29').Value = '38.50'
$ws.Range('E29').Value = '  -4.42%  '
$ws.Range('D30').Value = '22.38'
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').Value = '173.33'
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('D32').Value = '3.16'
$ws.Range('E32').Value = '  -3.33%  '
$ws.Range('D33').Value = '0.0901'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').Value = '5.74'
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').Value = '5.03'
$ws.Range('E35').Value = '  +2.97%  '
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('D37').Value = '4.12'
$ws.Range('E37').Value = '  -9.08%  '
$ws.Range('D38').Value = '0.0374'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  -2.27%  '
$ws.Range('D40').Value = '2.64'
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('D41').Value = '1.49'
$ws.Range('E41').Value = '  +10.80%  '
$ws.Range('D42').Value = '71.42'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').Value = '0.232'
$ws.Range('E43').Value = '  -4.28%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '12.32'
$ws.Range('E45').Value = '  -8.57%  '
$ws.Range('D46').Value = '5.73'
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('D47').Value = '108.90'
$ws.Range('E47').Value = '  +1.82%  '
$ws.Range('D48').Value = '9.00'
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('E49').Value = '  -6.03%  '
$ws.Range('D50').Value = '0.0988'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').Value = '70.37'
$ws.Range('E51').Value = '  +0.46%  '
